$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.169.42'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '2.010.49'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '258.15'
$ws.Range("E5").Value = '  +4.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.619'
$ws.Range("E6").Value = '  -1.31%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.17'
$ws.Range("E8").Value = '  -7.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.379'
$ws.Range("E9").Value = '  -2.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0768'
$ws.Range("E10").Value = '  -5.04%  '
$ws.Range("E11").Value = '  -2.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.32'
$ws.Range("E12").Value = '  -4.98%  '
$ws.Range("D13").Value = '2.309.78'
$ws.Range("E13").Value = '  -0.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.49'
$ws.Range("E14").Value = '  -1.91%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.798'
$ws.Range("E15").Value = '  -6.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.21'
$ws.Range("D17").Value = '2.020.73'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").Value = '37.012.30'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.97'
$ws.Range("E19").Value = '  +0.82%  '
$ws.Range("D20").Value = '0.0₃0834'
$ws.Range("E20").Value = '  -3.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '234.14'
$ws.Range("E21").Value = '  +1.67%  '
$ws.Range("E22").Value = '  -2.53%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.56'
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.05'
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.96'
$ws.Range("E27").Value = '  -4.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.52'
$ws.Range("E28").Value = '  -1.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.126'
$ws.Range("E29").Value = '  -9.06%  '
$ws.Range("E30").Value = '  -3.58%  '
$ws.Range("E31").Value = '  -2.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.59'
$ws.Range("E32").Value = '  -3.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0641'
$ws.Range("E33").Value = '  -4.92%  '
$ws.Range("E34").Value = '  -1.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.36'
$ws.Range("E35").Value = '  -5.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.49'
$ws.Range("E36").Value = '  -3.91%  '
$ws.Range("E37").Value = '  +0.91%  '
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.53'
$ws.Range("E39").Value = '  +3.66%  '
$ws.Range("E40").Value = '  +0.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.18'
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("D42").Value = '1.440.92'
$ws.Range("E42").Value = '  +4.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '89.40'
$ws.Range("E45").Value = '  -2.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.57'
$ws.Range("E46").Value = '  -8.06%  '
$ws.Range("E47").Value = '  -3.11%  '
$ws.Range("E48").Value = '  +2.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.98'
$ws.Range("E49").Value = '  -6.33%  '
$ws.Range("D50").Value = '2.200.96'
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("E51").Value = '  -7.70%  '

# Row 43/44 swap (Cronos <-> VeChain) with updated values
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0210'
$ws.Range("E43").Value = '  -2.52%  '

$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0921'
$ws.Range("E44").Value = '  -5.55%  '
